$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46011
$ws.Range("B2").Value = 92.98
$ws.Range("C2").Value = 83.51000000000001
$ws.Range("D2").Value = 77.63
$ws.Range("E2").Value = 68.25
$ws.Range("F2").Value = 67.59
$ws.Range("G2").Value = 67.06
$ws.Range("H2").Value = 73.2
$ws.Range("I2").Value = 75.55
$ws.Range("J2").Value = 79.25
$ws.Range("K2").Value = 80.84
$ws.Range("L2").Value = 75.03
$ws.Range("M2").Value = 65.08
$ws.Range("N2").Value = 59.85
$ws.Range("O2").Value = 49.87
$ws.Range("P2").Value = 52.52
$ws.Range("Q2").Value = 54.92
$ws.Range("R2").Value = 63.72
$ws.Range("S2").Value = 74.58
$ws.Range("T2").Value = 76.44
$ws.Range("U2").Value = 77.54000000000001
$ws.Range("V2").Value = 73.28
$ws.Range("W2").Value = 69.41
$ws.Range("X2").Value = 73.17
$ws.Range("Y2").Value = 67.2
$ws.Range("Z2").Value = 70.77
$ws.Range("AA2").Value = "0h-4h"
$ws.Range("AB2").Value = 80.59
$ws.Range("AC2").Value = "0h-2h"
$ws.Range("AD2").Value = 88.23999999999999
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 80.04000000000001
$ws.Range("AG2").Value = "3h-23h"
